# Applies the weekly update to the Hortaliza / Vega Monumental Concepción - Albahaca sheet.
# The rows effectively got re-sorted/updated: for each data row (2,3,4,6-15) the
# columns D (Fecha), J (Volumen), K (Precio mínimo), L (Precio máximo),
# M (Precio promedio ponderado), O (Origen) and P (Precio $/Kg) take on new values.
# Row 5 is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map: row number -> @(D, J, K, L, M, O, P)
$updates = @{
    2  = @(44643, 90,  2800, 3000, 2911, "Región Metropolitana",   485)
    3  = @(44650, 130, 3000, 3500, 3308, "Región Metropolitana",   551)
    4  = @(44659, 90,  2500, 3000, 2722, "Región Metropolitana",   454)
    6  = @(44957, 70,  1500, 2000, 1857, "Región Metropolitana",   310)
    7  = @(44685, 150, 3000, 3500, 3267, "Región Metropolitana",   544)
    8  = @(44876, 80,  6500, 7000, 6812, "Región Metropolitana",  1135)
    9  = @(44672, 140, 3000, 3500, 3286, "Región Metropolitana",   548)
    10 = @(44671, 150, 3500, 4000, 3733, "Región Metropolitana",   622)
    11 = @(44631, 110, 3000, 3500, 3273, "Provincia de Chacabuco", 546)
    12 = @(44644, 140, 2500, 3000, 2786, "Provincia de Chacabuco", 464)
    13 = @(44637, 170, 2800, 3000, 2906, "Región Metropolitana",   484)
    14 = @(44658, 180, 2500, 3000, 2778, "Región Metropolitana",   463)
    15 = @(44987, 130, 4500, 5000, 4692, "Región Metropolitana",   782)
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]

    $ws.Range("D$row").Value = $vals[0]
    $ws.Range("J$row").Value = $vals[1]
    $ws.Range("K$row").Value = $vals[2]
    $ws.Range("L$row").Value = $vals[3]
    $ws.Range("M$row").Value = $vals[4]
    $ws.Range("O$row").Value = $vals[5]
    $ws.Range("P$row").Value = $vals[6]
}
